$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two JSON request bodies and the path param value
$ws.Range("B3").Value = '{   "firstName": "stefan",   "lastName": "maclory",   "phoneNumber": "00000099999",   "emailAddress": "jone.mac@marvel.com" }'
$ws.Range("B4").Value = '{   "firstName": "danis",   "lastName": "targerian",   "phoneNumber": "8888899999",   "emailAddress": "jone.mac@thripaty.com" }'
$ws.Range("D4").Value = "/3"

# Update the selection/view: select B3 (also resets the scrolled topLeftCell back to default)
$ws.Range("B3").Select()
